# Update "Riders" (column C) and "Average" (column D) values on the
# Ridership sheet for Nov 2016 — new Madigan bike hours.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$data = @(
    @(2, 260, 260),
    @(3, 270, 270),
    @(4, 266, 266),
    @(5, 234, 257.5),
    @(6, 96, 111.5),
    @(7, 72, 78),
    @(8, 257, 263),
    @(9, 209, 234.5),
    @(10, 232, 251),
    @(11, 221, 243.5),
    @(12, 116, 113),
    @(13, 66, 74),
    @(14, 306, 277.33),
    @(15, 212, 227),
    @(16, 264, 255.33),
    @(17, 261, 249.33),
    @(18, 246, 253.67),
    @(19, 127, 116.5),
    @(20, 86, 77),
    @(21, 233, 266.25),
    @(22, 209, 222.5),
    @(23, 200, 241.5),
    @(24, 172, 233.25),
    @(25, 131, 119.4),
    @(26, 96, 80.8),
    @(27, 246, 262.2),
    @(28, 277, 233.4),
    @(29, 20, 197.2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
